# Netherlands.xlsx data fix-up: row 3 (id 150) had an oldest-known date of
# 2007-02-18 sourced from "abacq date posted"; correct it to just the year
# 1974 (month/day unknown) with a new source label "desc other".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# oldest_known_year (N3): 2007 -> 1974
$ws.Range("N3").Value = 1974

# oldest_known_month / oldest_known_day (O3, P3): no longer known, clear them
$ws.Range("O3").ClearContents()
$ws.Range("P3").ClearContents()

# oldest_known_source (Q3): "abacq date posted" -> new string "desc other"
$ws.Range("Q3").Value = "desc other"

# Reflect the author's last selection in the saved view
$ws.Range("A3").Select()
